$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.082.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.45%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.435.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.77%  "
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.89%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.54"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +12.02%  "
# Row 7
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.605"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.92%  "
# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.427.09"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.68%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.681"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.94%  "
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +10.11%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.58"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.42%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.142"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.39%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.970.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.17%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.92%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.88"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.21%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.489.04"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.73%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.988.33"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.68%  "
# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.03"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.45%  "
# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.46%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000138"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +15.01%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.27"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.54%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.63"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.42%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.06"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.30%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "310.36"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.11%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.19"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.78%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.72"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +9.60%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.85"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.12%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.45"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.90%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.54"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.82%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.175"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.27%  "
# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.42%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "43.38"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.33%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.73"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.22%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.17%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.16%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0484"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.56%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.44"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.68%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.18%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.47"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.33%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.97"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.37%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.126"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.65%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.51"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.15%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.73%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.05%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.285"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.48%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.91"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.34%  "
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.93%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.91"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.47%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.772.63"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.49%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.151.55"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.30%  "

Write-Host "Applied cryptos update"